$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": add a new day column CE (04-sep) with header + 24 values
# ---------------------------------------------------------------------------
$wsSpot = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting of the previous header cell (CD1) onto the new header
# cell (CE1) so it keeps the bold / centered / bordered look of the other
# header cells, then set its text.
$wsSpot.Range("CD1").Copy()
$wsSpot.Range("CE1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsSpot.Range("CE1").Value = "04-sep"

$spotValues = @(19.84, 14.81, 14.3, 5.96, 5.99, 11.49, 42.15, 52.04, 55.89, 49.6, 18.07, 5.79, 3.2, 1.72, 1.6, 1.55, 3.7, 7.01, 28.19, 62.19, 102.5, 80.74, 84.8, 62.18)
for ($i = 0; $i -lt $spotValues.Length; $i++) {
    $row = $i + 2
    $wsSpot.Cells.Item($row, 83).Value = $spotValues[$i]
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append a new row 80 (2025-09-02, 30.5)
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A80").NumberFormat = "@"
$wsGaz.Range("A80").Value = "2025-09-02"
$wsGaz.Range("A80").Style = "Normal"
$wsGaz.Range("B80").Value = 30.5

# ---------------------------------------------------------------------------
# Sheet "CO2": append a new row 80 (2025-09-02, 73.25)
# ---------------------------------------------------------------------------
$wsCO2 = $wb.Worksheets.Item("CO2")
$wsCO2.Range("A80").NumberFormat = "@"
$wsCO2.Range("A80").Value = "2025-09-02"
$wsCO2.Range("A80").Style = "Normal"
$wsCO2.Range("B80").Value = 73.25
